# Purchase temp related issue solution
# - Remove 9 obsolete/duplicate purchase line rows (old rows 16-24)
# - Relabel the "Expiry Date" header to "EXPIRE Date"
# - Restore the normal frozen-pane/selection view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterSheet")

# Delete the 9 rows that were removed from the purchase list
# (AG10, AG23, B1, AH26-dup, HD67, BG9, A92, W73, A86-dup)
$ws.Rows("16:24").Delete()

# Rename header label in D1
$ws.Range("D1").Value = "EXPIRE Date"

# Restore view (frozen pane anchored back at row 2, selection on E8)
$ws.Range("E8").Select()
